# Apply the "export client pertama push" edit to the workbook.
#
# Summary of changes:
#   1. Rename the only worksheet from "Sheet1" to "Bank KB Bukopin".
#   2. Remove all the stale #REF! defined names left over in the workbook.
#   3. Zoom the sheet view from 70% to 90%.
#   4. Move the active selection on the sheet from A7:XFD7 to the single
#      cell H14.
#   5. Give row 6 an explicit height of 15.75 (matching row 5's height).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet.
$ws.Name = "Bank KB Bukopin"

# 2. Remove every defined name in the workbook (asli02, ASLI08, ASLI1,
#    ASLI3, DATA1, DATABASPG, DOAB, Frezee_Employee, mei,
#    Tagihan_BPJS_Kesehatan, tdk, TES - including their sheet-scoped
#    duplicates). Always delete item 1 since the collection re-indexes
#    after every removal.
while ($wb.Names.Count -gt 0) {
    $wb.Names.Item(1).Delete()
}

# 3. Update the zoom level for the sheet's window from 70% to 90%.
$excel.ActiveWindow.Zoom = 90

# 4. Change the active selection to H14 (previously the whole row A7:XFD7).
$ws.Range("H14").Select()

# 5. Set row 6's height to 15.75pt, matching row 5.
$ws.Rows.Item(6).RowHeight = 15.75
